# Updates the crypto price ("D" column) and 1h volume % change ("E" column)
# values for rows 2-51 on the active worksheet, matching the latest scrape
# from the "Updated cryptos list ... with GitHub Actions" commit.
#
# Cells whose new text would otherwise be auto-recognized by Excel as a
# number (losing trailing zeros / changing cell type) are temporarily
# switched to Text format, written, then restored to their original style
# so the resulting cell keeps being a plain text value (matching the
# original workbook's inline/shared string cells) without any unintended
# style/number-format changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.106.17"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.791.64"
$ws.Range("E3").Value = "  -1.41%  "
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.37%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.56"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("E6").Value = "  +0.34%  "
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5223"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +2.68%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3794"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  -4.06%  "
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07952"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  -3.44%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.38"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -1.28%  "
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +0.35%  "
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.245"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -0.57%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.46"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "1.793.01"
$ws.Range("E15").Value = "  -1.52%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.279"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -3.11%  "
$ws.Range("E17").Value = "  -1.00%  "
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001088"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -4.48%  "
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06565"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("E20").Value = "  +0.34%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.27"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -2.28%  "
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.960"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "28.157.41"
$ws.Range("E23").Value = "  -1.05%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.09"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -2.21%  "
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.253"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -0.62%  "
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.08"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +2.97%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.39"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -3.55%  "
$ws.Range("D28").Value = "1.997.57"
$ws.Range("E28").Value = "  -1.43%  "
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.331"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -3.02%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.82"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -2.38%  "
$ws.Range("E31").Value = "  -0.77%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.054"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("E33").Value = "  +1.13%  "
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.518"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  +2.05%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.17"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +8.05%  "
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02307"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  -3.39%  "
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.051"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -3.11%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.580"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("E41").Value = "  -2.01%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.159"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -1.41%  "
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.372"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -2.29%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.19"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -2.29%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.763"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +0.59%  "
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5907"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +0.20%  "
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.51"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +1.17%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.211"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +2.01%  "
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.914"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -3.04%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06779"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -1.60%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.47"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -2.18%  "
